$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# 1. Normalize the D2:D63 "Y" column style (drop the redundant fill flag
#    that Excel had been carrying on those cells) by repainting the format
#    from a cell that already uses the clean style.
$ws.Range("A2").Copy()
$ws.Range("D2:D63").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. OPQA-378 test case (row 59) moved from PASS to SKIP.
$ws.Range("E59").Value2 = "SKIP"

# 3. Append the new profile test cases (rows 64-68).
$ws.Range("A64").Value2 = "TestCase_B63"
$ws.Range("B64").Value2 = "OPQA-258"
$ws.Range("C64").Value2 = "Verify that no search results get displayed if search engine doesn't interpret the query and that a proper message gets displayed regarding that"
$ws.Range("D64").Value2 = "Y"
$ws.Range("E64").Value2 = "SKIP"

$ws.Range("A65").Value2 = "TestCase_B64"
$ws.Range("B65").Value2 = "OPQA-557"
$ws.Range("C65").Value2 = "Verify that the searched keyword present in the search text box doesn't change if any other content type is selected in the search drop down"
$ws.Range("D65").Value2 = "Y"
$ws.Range("E65").Value2 = "SKIP"

$ws.Range("A66").Value2 = "TestCase_B65"
$ws.Range("B66").Value2 = "OPQA-386"
$ws.Range("C66").Value2 = "Verify that the searched keyword doesn't change in the search text box if any other content type is selected in the left navigation pane"
$ws.Range("D66").Value2 = "Y"
$ws.Range("E66").Value2 = "SKIP"

$ws.Range("A67").Value2 = "TestCase_B66"
$ws.Range("B67").Value2 = "OPQA-387"
$ws.Range("C67").Value2 = "Verify that counts of search results of all the content types should get displayed irrespective of the content type chosen for searching"
$ws.Range("D67").Value2 = "Y"
$ws.Range("E67").Value2 = "SKIP"

$ws.Range("A68").Value2 = "TestCase_B67"
$ws.Range("B68").Value2 = "OPQA-263"
$ws.Range("C68").Value2 = "Verify that ALL search results count is equal to the count of search results of other content types(ARTICLES+PATENTS+POSTS+PEOPLE)"
$ws.Range("D68").Value2 = "Y"
$ws.Range("E68").Value2 = "PASS"

# 4. Apply the same border/format to the new rows as the rest of the table.
$ws.Range("A63:E63").Copy()
$ws.Range("A64:E68").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. Restore values that PasteSpecial (formats) doesn't touch -- they are
#    already set above, so nothing else needed here. Re-assert the
#    selection the author left on the sheet.
$ws.Range("D10").Select()
